$d = $word.ActiveDocument

$replacements = @(
    @{old = "56×91="; new = "26×97="},
    @{old = "93×45="; new = "17×60="},
    @{old = "85×28="; new = "18×35="},
    @{old = "63×96="; new = "30×25="},
    @{old = "71×32="; new = "78×54="},
    @{old = "45×87="; new = "40×31="},
    @{old = "29×49="; new = "80×60="},
    @{old = "55×20="; new = "85×33="},
    @{old = "65×34="; new = "92×68="},
    @{old = "15×88="; new = "75×11="},
    @{old = "52×53="; new = "15×29="},
    @{old = "26×49="; new = "42×82="},
    @{old = "41×22="; new = "67×42="},
    @{old = "75×36="; new = "23×74="},
    @{old = "56×87="; new = "21×32="},
    @{old = "58×96="; new = "72×93="},
    @{old = "13×19="; new = "37×43="},
    @{old = "74×90="; new = "56×73="},
    @{old = "71×26="; new = "13×55="},
    @{old = "37×44="; new = "21×62="},
    @{old = "93×91="; new = "60×93="},
    @{old = "39×45="; new = "21×89="},
    @{old = "98×34="; new = "38×48="},
    @{old = "70×12="; new = "53×38="},
    @{old = "27×48="; new = "93×60="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $r.new, 2)
}
